# Add a new day's data to the top of the table (row 2), pushing existing
# rows down by one, matching how the sheet is updated each day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the first data row). This
# shifts rows 2-38 down to 3-39, and Excel automatically adjusts the
# relative formulas/shared-formula ranges and row references.
$ws.Rows.Item(2).Insert()

# Preserve the date number format on the new A2 cell (columns B:K already
# inherit their number format from the column styles).
$ws.Range("A2").NumberFormat = "d-mmm"

# Fill in the new day's values.
$ws.Range("A2").Value = 44284
$ws.Range("B2").Value = 0.1296
$ws.Range("C2").Value = 0.1618
$ws.Range("D2").Value = 0.0993
$ws.Range("E2").Value = 0.115
$ws.Range("F2").Value = 0.0991
$ws.Range("G2").Value = 0.0566
$ws.Range("H2").Value = 0.0582
$ws.Range("I2").Value = 0.0892
$ws.Range("J2").Value = 44284
$ws.Range("K2").Value = 30

# Match the saved selection/active cell recorded in the workbook.
[void]$ws.Range("M12").Select()
